$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Header row (row 1): two new leading "Depth from/to [m]" columns
# replace the old "z from/to [ft]" columns (A1/B1 keep their bold
# bordered style), and the last header (H1) is renamed from
# "Subm unit weight [kN/m3]" to "Total unit weight [kN/m3]".
# ------------------------------------------------------------------
$ws.Range("A1").Value = "Depth from [m]"
$ws.Range("B1").Value = "Depth to [m]"
$ws.Range("H1").Value = "Total unit weight [kN/m3]"

# The remaining header cells (C1:H1) get a new bold, border-less style.
$ws.Range("C1:H1").Font.Bold = $true

# ------------------------------------------------------------------
# Relative-density text cells - values unchanged, just re-written so
# the shared-string table is rebuilt the same way Excel rebuilt it.
# ------------------------------------------------------------------
$ws.Range("D2").Value = "Loose"
$ws.Range("D4").Value = "Medium dense"
$ws.Range("D5").Value = "Dense"

# ------------------------------------------------------------------
# Unit weight column now reports "Total" instead of "Submerged"
# values (+10 kN/m3 across the board).
# ------------------------------------------------------------------
$ws.Range("H2").Value = 19
$ws.Range("H3").Value = 18
$ws.Range("H4").Value = 19
$ws.Range("H5").Value = 20

# ------------------------------------------------------------------
# Column widths - re-sized/added to fit the new headers (best-fit
# character widths as computed by Excel for the new content).
# ------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 12.25
$ws.Columns.Item(2).ColumnWidth = 10.25
$ws.Columns.Item(3).ColumnWidth = 6.92
$ws.Columns.Item(4).ColumnWidth = 12.6
$ws.Columns.Item(5).ColumnWidth = 11.25
$ws.Columns.Item(6).ColumnWidth = 9.25
$ws.Columns.Item(7).ColumnWidth = 7.09
$ws.Columns.Item(8).ColumnWidth = 20.59

# ------------------------------------------------------------------
# Selection moved from M13 to H13.
# ------------------------------------------------------------------
$ws.Range("H13").Select() | Out-Null

Write-Output "done"
